$d = $word.ActiveDocument

# Fix the table caption paragraph text
$d.Content.Find.Execute(
    "(#tab:unnamed-chunk-13)Descripte Statistics", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "(#tab:unnamed-chunk-3)Descriptive Statistics", 2)

# Fix the table's Title (maps to w:tblCaption)
$tbl = $d.Tables.Item(1)
$tbl.Title = "(#tab:unnamed-chunk-3)Descriptive Statistics"

# Fix the statistic values
$d.Content.Find.Execute("0.008", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.002", 2)

$d.Content.Find.Execute("0.195 (0.399)", $true, $false, $false, $false,
                         $false, $true, 1, $false, "0.244 (0.432)", 2)

$d.Content.Find.Execute("0.408 (0.497)", $true, $false, $false, $false,
                         $false, $true, 1, $false, "0.510 (0.505)", 2)

$d.Content.Find.Execute("0.275 (0.448)", $true, $false, $false, $false,
                         $false, $true, 1, $false, "0.344 (0.477)", 2)
